$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.380.14'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").Value = '1.868.51'
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''330.51'
$ws.Range("E5").Value = '  -2.24%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").Value = '''0.4594'
$ws.Range("E7").Value = '  -2.25%  '

$ws.Range("D8").Value = '''0.4006'
$ws.Range("E8").Value = '  +1.90%  '

$ws.Range("D9").Value = '''47.70'
$ws.Range("E9").Value = '  +1.27%  '

$ws.Range("D10").Value = '''0.07843'
$ws.Range("E10").Value = '  -1.79%  '

$ws.Range("D11").Value = '''0.9838'
$ws.Range("E11").Value = '  -1.90%  '

$ws.Range("D12").Value = '''21.27'
$ws.Range("E12").Value = '  -2.29%  '

$ws.Range("D13").Value = '1.882.26'
$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").Value = '''5.847'
$ws.Range("E14").Value = '  -2.41%  '

$ws.Range("D15").Value = '''6.988'
$ws.Range("E15").Value = '  -3.97%  '

$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").Value = '''88.20'
$ws.Range("E17").Value = '  -3.19%  '

$ws.Range("D18").Value = '''0.06526'
$ws.Range("E18").Value = '  -0.86%  '

$ws.Range("D19").Value = '''0.00001017'
$ws.Range("E19").Value = '  -2.53%  '

$ws.Range("D20").Value = '''17.16'
$ws.Range("E20").Value = '  -2.89%  '

$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").Value = '28.369.09'
$ws.Range("E22").Value = '  +0.26%  '

$ws.Range("E23").Value = '  -2.31%  '

$ws.Range("D24").Value = '''10.84'
$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").Value = '''2.250'
$ws.Range("E25").Value = '  -1.78%  '

$ws.Range("D26").Value = '2.098.19'
$ws.Range("E26").Value = '  -0.70%  '

$ws.Range("D27").Value = '''157.46'
$ws.Range("E27").Value = '  -1.20%  '

$ws.Range("E28").Value = '  -2.81%  '

$ws.Range("E29").Value = '  -4.80%  '

$ws.Range("D30").Value = '''5.278'
$ws.Range("E30").Value = '  -4.06%  '

$ws.Range("D31").Value = '''117.23'

$ws.Range("D32").Value = '''0.9534'
$ws.Range("E32").Value = '  -2.73%  '

$ws.Range("D33").Value = '''0.09311'
$ws.Range("E33").Value = '  -1.87%  '

$ws.Range("D34").Value = '''3.591'
$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("D35").Value = '''1.384'
$ws.Range("E35").Value = '  +0.14%  '

$ws.Range("D36").Value = '''5.221'
$ws.Range("E36").Value = '  -2.55%  '

$ws.Range("D37").Value = '''0.06018'
$ws.Range("E37").Value = '  -1.21%  '

$ws.Range("D38").Value = '''0.02199'
$ws.Range("E38").Value = '  -3.04%  '

$ws.Range("D39").Value = '''8.265'
$ws.Range("E39").Value = '  -2.13%  '

$ws.Range("D40").Value = '''1.167'
$ws.Range("E40").Value = '  -0.82%  '

$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("D42").Value = '''0.5745'
$ws.Range("E42").Value = '  -3.68%  '

$ws.Range("D43").Value = '''0.1805'
$ws.Range("E43").Value = '  -3.84%  '

$ws.Range("D44").Value = '''9.987'
$ws.Range("E44").Value = '  -3.79%  '

$ws.Range("D45").Value = '''1.261'
$ws.Range("E45").Value = '  -3.13%  '

$ws.Range("D46").Value = '''2.273'
$ws.Range("E46").Value = '  +12.67%  '

$ws.Range("E47").Value = '  -3.52%  '

$ws.Range("D48").Value = '''11.80'
$ws.Range("E48").Value = '  -2.96%  '

$ws.Range("D49").Value = '''0.07144'
$ws.Range("E49").Value = '  +3.57%  '

$ws.Range("D50").Value = '''1.881'
$ws.Range("E50").Value = '  -4.50%  '

$ws.Range("D51").Value = '''110.19'
$ws.Range("E51").Value = '  -0.52%  '
